# Update column F ("dSF") values for the specified rows.
# This corresponds to the commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    4  = -2
    6  = 3
    12 = 3
    13 = -5
    17 = 1
    20 = 3
    26 = -2
    27 = 2
    28 = 1
    31 = 3
    39 = 3
    46 = 2
    50 = -3
    51 = 2
    55 = -2
    63 = 1
    65 = 0
    70 = 2
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
